$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.018.08'
$ws.Range("E2").Value = '  +2.22%  '
$ws.Range("D3").Value = '1.675.77'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = '''329.86'
$ws.Range("E5").Value = '  +7.32%  '
$ws.Range("D6").Value = '''0.9993'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").Value = '''0.3657'
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("D8").Value = '''47.38'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").Value = '''0.3239'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '''1.147'
$ws.Range("E10").Value = '  +2.61%  '
$ws.Range("D11").Value = '''0.07169'
$ws.Range("E11").Value = '  +2.40%  '
$ws.Range("D12").Value = '''0.9996'
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '''6.090'
$ws.Range("E13").Value = '  +3.56%  '
$ws.Range("D14").Value = '''19.68'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '1.670.33'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").Value = '''0.00001049'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '''0.06537'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '''0.9990'
$ws.Range("D20").Value = '''78.93'
$ws.Range("E20").Value = '  +3.46%  '
$ws.Range("D21").Value = '''15.84'
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("D22").Value = '''5.913'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '''12.93'
$ws.Range("E23").Value = '  +2.84%  '
$ws.Range("D24").Value = '25.005.00'
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("D25").Value = '''2.437'
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("D26").Value = '''2.378'
$ws.Range("E26").Value = '  +2.85%  '
$ws.Range("D27").Value = '''149.12'
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("D28").Value = '''18.74'
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").Value = '1.857.46'
$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("D30").Value = '''125.91'
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("D31").Value = '''1.199'
$ws.Range("E31").Value = '  +1.76%  '
$ws.Range("D32").Value = '''4.089'
$ws.Range("E32").Value = '  +2.61%  '
$ws.Range("D33").Value = '''5.796'
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("D35").Value = '''1.670'
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").Value = '''12.34'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").Value = '''5.157'
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("D38").Value = '''0.06079'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").Value = '''1.233'
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("D40").Value = '''0.2091'
$ws.Range("E40").Value = '  +1.81%  '
$ws.Range("D41").Value = '''0.02229'
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("D42").Value = '''8.265'
$ws.Range("E42").Value = '  +1.09%  '
$ws.Range("D43").Value = '''0.9985'
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").Value = '''0.5968'
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("D45").Value = '''13.58'
$ws.Range("E45").Value = '  +7.95%  '
$ws.Range("D46").Value = '''3.828'
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("D47").Value = '''0.5736'
$ws.Range("E47").Value = '  +2.73%  '
$ws.Range("D48").Value = '''124.04'
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("D49").Value = '''1.967'
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").Value = '''0.07014'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").Value = '''1.189'
$ws.Range("E51").Value = '  +3.19%  '
